$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

# New header for column S
$ws.Range("S1").Value = "strNameLen"

# Update row 2 values (data row)
$ws.Range("A2").Value = "testPlate.dxf"
$ws.Range("B2").Value = "NTNU_KeyChain_Template_v2.dxf"
$ws.Range("C2").Value = "ntnu_logo_svart.png"

# D2, I2, M2, R2 use the "Text" number format style already present in style index 1.
# Apply that style (NumberFormat "@") before setting the values so they are stored as text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.5050"

$ws.Range("E2").Value = "18.7732"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "29.5"

$ws.Range("K2").Value = "JAGAMcGCHE-J"
$ws.Range("L2").Value = "NTNU-DIN"

$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "12.5"

$ws.Range("N2").Value = 14

$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "2.5"

$ws.Range("S2").Value = 8

# Update window view settings
$wb.Windows.Item(1).WindowState = -4143
$excel.ActiveWindow.WindowState = -4143

# Set sheet view: topLeftCell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("R2").Select()
